$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 9 so the existing totals row (old r9) and
# the footer row (old r10) shift down to r11/r12, matching the diff.
$ws.Rows("9:10").Insert()

# Copy the formatting (styles/borders/number formats/merges) from the two
# product rows above (7 and 8) down onto the freshly inserted rows 9 and 10
# so the new product rows look the same as the existing ones.
$ws.Range("A7:Q7").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Range("A8:Q8").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)

# Row heights for the two new rows match rows 7 / 8 respectively.
$ws.Rows("9").RowHeight = 25.5
$ws.Rows("10").RowHeight = 24.75

# Row 9 - third product line: فرشه اسنان POWER GOLD كبار
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "فرشه اسنان POWER GOLD كبار"
$ws.Range("H9").Value = "0:0"
$ws.Range("L9").Value = "0"
$ws.Range("N9").Value = "15.00"
$ws.Range("P9").Value = "15.0000"
$ws.Range("Q9").Value = "1:0"

# Row 10 - fourth product line: معجون سيجنال 25 مل
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "معجون سيجنال 25 مل"
$ws.Range("H10").Value = "6:0"
$ws.Range("L10").Value = "0"
$ws.Range("N10").Value = "20.00"
$ws.Range("P10").Value = "20.0000"
$ws.Range("Q10").Value = "1:0"

# Row 11 (was row 9) - grand total updates to include the two new lines.
$ws.Range("P11").Value = 66.83

# Row 12 (was row 10) - footer: refresh the generated timestamp.
$ws.Range("A12").Value = "Friday, 18 July, 2025 4:26 PM"
